$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 11: Test Type changes from "Functional" to "Login"
$ws.Range("C11").Value = "Login"

# Add new row 12 for the "Navigate_AllMenuItems" test case
# (order of cell writes controls the order new shared strings are interned)
$ws.Range("C12").Value = "Login"
$ws.Range("D12").Value = "Verif successful navigation to different menu items available in kirana Bazaar application."
$ws.Range("B12").Value = "Navigate_AllMenuItems"
$ws.Range("A12").Value = 11

# Match formatting/style used by the other "Test Description" cells (wrap text)
$ws.Range("D12").WrapText = $true

# Update the selected/active cell to match the recorded view state
$ws.Range("D7").Select()
